$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "TextBox 7") {
        $shp.Delete()
        break
    }
}
